# Refresh market-price derived columns (H-N) across the Leve profit sheets.
# Mirrors the scheduled runner's scrape-and-recompute pass: plain numeric
# values only (no formulas in these tables), a few rows also gain/lose a
# trailing HQ-profit cell (N) as that branch's computed value's availability
# changes from run to run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1388.7142
$ws.Range("I12").Value = 911
$ws.Range("K12").Value = 911
$ws.Range("M12").Value = -741
$ws.Range("H132").Value = 2670.85
$ws.Range("I132").Value = 2682.9722
$ws.Range("J132").Value = 2561.75
$ws.Range("K132").Value = 8048.9166
$ws.Range("L132").Value = 7685.25
$ws.Range("M132").Value = -5518.9166
$ws.Range("N132").Value = -12745.25
$ws.Range("H137").Value = 2063.359
$ws.Range("I137").Value = 1932.88
$ws.Range("J137").Value = 2296.3572
$ws.Range("K137").Value = 5798.64
$ws.Range("L137").Value = 6889.071599999999
$ws.Range("M137").Value = -3248.64
$ws.Range("N137").Value = -11989.0716
$ws.Range("H138").Value = 2324.9756
$ws.Range("I138").Value = 1749.5652
$ws.Range("K138").Value = 5248.6956
$ws.Range("M138").Value = -108.6956
$ws.Range("H141").Value = 1900.7142
$ws.Range("I141").Value = 1950.1621
$ws.Range("K141").Value = 5850.4863
$ws.Range("M141").Value = -670.4863000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 34999
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 34999
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 34999
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -35457
$ws.Range("H37").Value = 29999.5
$ws.Range("J37").Value = 49999
$ws.Range("L37").Value = 49999
$ws.Range("N37").Value = -50545

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5435365
$ws.Range("I94").Value = 7143329.5
$ws.Range("J94").Value = 931.63635
$ws.Range("K94").Value = 7143329.5
$ws.Range("L94").Value = 931.63635
$ws.Range("M94").Value = -7142878.5
$ws.Range("N94").Value = -1833.63635
$ws.Range("H134").Value = 3792.9456
$ws.Range("I134").Value = 3524.3171
$ws.Range("J134").Value = 4579.643
$ws.Range("K134").Value = 10572.9513
$ws.Range("L134").Value = 13738.929
$ws.Range("M134").Value = -8037.951300000001
$ws.Range("N134").Value = -18808.929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1013
$ws.Range("I22").Value = 907.25
$ws.Range("J22").Value = 1062.7646
$ws.Range("K22").Value = 907.25
$ws.Range("L22").Value = 1062.7646
$ws.Range("M22").Value = -557.25
$ws.Range("N22").Value = -1762.7646
$ws.Range("H33").Value = 705.1667
$ws.Range("I33").Value = 705.1667
$ws.Range("K33").Value = 705.1667
$ws.Range("M33").Value = -326.1667
$ws.Range("H39").Value = 15493.75
$ws.Range("I39").Value = 7658.3335
$ws.Range("J39").Value = 39000
$ws.Range("K39").Value = 7658.3335
$ws.Range("L39").Value = 39000
$ws.Range("M39").Value = -7267.3335
$ws.Range("N39").Value = -39782
$ws.Range("H49").Value = 15493.75
$ws.Range("I49").Value = 7658.3335
$ws.Range("J49").Value = 39000
$ws.Range("K49").Value = 7658.3335
$ws.Range("L49").Value = 39000
$ws.Range("M49").Value = -7476.3335
$ws.Range("N49").Value = -39364
$ws.Range("H58").Value = 2491.4849
$ws.Range("I58").Value = 2715.92
$ws.Range("K58").Value = 2715.92
$ws.Range("M58").Value = -2512.92
$ws.Range("H60").Value = 38981
$ws.Range("J60").Value = 38981
$ws.Range("L60").Value = 38981
$ws.Range("N60").Value = -40003
$ws.Range("H99").Value = 3269.3684
$ws.Range("J99").Value = 2662.6667
$ws.Range("L99").Value = 2662.6667
$ws.Range("N99").Value = -5658.6667
$ws.Range("H126").Value = 3269.3684
$ws.Range("J126").Value = 2662.6667
$ws.Range("L126").Value = 7988.000100000001
$ws.Range("N126").Value = -12928.0001
$ws.Range("H136").Value = 2491.4849
$ws.Range("I136").Value = 2715.92
$ws.Range("K136").Value = 8147.76
$ws.Range("M136").Value = -5597.76

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 21514.5
$ws.Range("J105").Value = 21514.5
$ws.Range("L105").Value = 64543.5
$ws.Range("N105").Value = -69785.5
$ws.Range("H107").Value = 2219.5264
$ws.Range("J107").Value = 2380.8696
$ws.Range("L107").Value = 7142.6088
$ws.Range("N107").Value = -10982.6088
$ws.Range("H129").Value = 5666.1665
$ws.Range("I129").Value = 1356
$ws.Range("J129").Value = 8744.857
$ws.Range("K129").Value = 4068
$ws.Range("L129").Value = 26234.571
$ws.Range("M129").Value = 932
$ws.Range("N129").Value = -36234.571
$ws.Range("H132").Value = 1965.5555
$ws.Range("I132").Value = 1913
$ws.Range("K132").Value = 17217
$ws.Range("M132").Value = -14687
$ws.Range("H141").Value = 6577.1177
$ws.Range("J141").Value = 14999.25
$ws.Range("L141").Value = 44997.75
$ws.Range("N141").Value = -55357.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4445261
$ws.Range("I122").Value = 5922348.5
$ws.Range("K122").Value = 17767045.5
$ws.Range("M122").Value = -17764595.5
$ws.Range("H132").Value = 3252.7114
$ws.Range("I132").Value = 3190.3948
$ws.Range("J132").Value = 3421.8572
$ws.Range("K132").Value = 9571.1844
$ws.Range("L132").Value = 10265.5716
$ws.Range("M132").Value = -7041.1844
$ws.Range("N132").Value = -15325.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3336.7693
$ws.Range("I7").Value = 3190.1
$ws.Range("J7").Value = 3825.6667
$ws.Range("K7").Value = 3190.1
$ws.Range("L7").Value = 3825.6667
$ws.Range("M7").Value = -3078.1
$ws.Range("N7").Value = -4049.6667
$ws.Range("H40").Value = 7123.6
$ws.Range("I40").Value = 6900.8335
$ws.Range("J40").Value = 7457.75
$ws.Range("K40").Value = 6900.8335
$ws.Range("L40").Value = 7457.75
$ws.Range("M40").Value = -6764.8335
$ws.Range("N40").Value = -7729.75
$ws.Range("H61").Value = 8882.857
$ws.Range("J61").Value = 12260
$ws.Range("L61").Value = 12260
$ws.Range("N61").Value = -12664
$ws.Range("H113").Value = 8882.857
$ws.Range("J113").Value = 12260
$ws.Range("L113").Value = 12260
$ws.Range("N113").Value = -16600
$ws.Range("H126").Value = 3336.7693
$ws.Range("I126").Value = 3190.1
$ws.Range("J126").Value = 3825.6667
$ws.Range("K126").Value = 9570.299999999999
$ws.Range("L126").Value = 11477.0001
$ws.Range("M126").Value = -7100.299999999999
$ws.Range("N126").Value = -16417.0001
$ws.Range("H140").Value = 86536.60000000001
$ws.Range("J140").Value = 86536.60000000001
$ws.Range("L140").Value = 86536.60000000001
$ws.Range("N140").Value = -96896.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3403591
$ws.Range("I81").Value = 3403591
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6807182
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -6806121
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 3403591
$ws.Range("I84").Value = 3403591
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 34035910
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -34030606
$ws.Range("N84").ClearContents()
$ws.Range("H136").Value = 3634.7827
$ws.Range("I136").Value = 4032.5881
$ws.Range("J136").Value = 2507.6667
$ws.Range("K136").Value = 12097.7643
$ws.Range("L136").Value = 7523.000100000001
$ws.Range("M136").Value = -9547.764299999999
$ws.Range("N136").Value = -12623.0001

